# Improve formatting of test data for the oneWayFlightFilters test.
# D3 previously held a real date (1/1/2021) formatted with a date number
# format. It should instead hold the literal text "1-1-21", styled the
# same way as the existing text-date cell D2 (centered, quote-prefixed
# text so Excel doesn't reinterpret it as a date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy D2's formatting (font/alignment/quote-prefixed text style) onto D3,
# then set D3's value to the literal text "1-1-21" (leading apostrophe
# forces text entry so Excel does not convert it back into a date).
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D3").Value = "'1-1-21"

# Update the sheet's active selection from F6 to E8.
$ws.Range("E8").Select()
